$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns AD, AE, AF
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy the header style (bold, border, centered) from AC1 to the new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Fill in team record data for each player row (2-46)
for ($r = 2; $r -le 46; $r++) {
    $ws.Cells.Item($r, 30).Value = 77
    $ws.Cells.Item($r, 31).Value = 85
    $ws.Cells.Item($r, 32).Value = 0
}
